# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the period for the remaining worker (RAUL ENRIQUE GARRIDO SALAS) 2506 -> 2508
$ws.Range("E16").Value = "2508"

# Remove the second worker's data row (ESNAIDER MANJARREZ OLIVAREZ, doc 1002315868, period 2506)
$ws.Rows(17).Delete()

# Update summary totals now that only one worker/period remains
$ws.Range("E11").Value = 56940
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 1
